$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels for the new helper columns (J:M) on row 5 ---
$ws.Range("J5").Value = "C. Pagina"
$ws.Range("K5").Value = "C. Anillado"
$ws.Range("L5").Value = "Total"
$ws.Range("M5").Value = "Redondeado"

# --- Row 6 helper values + formulas ---
$ws.Range("J6").Value = 0.15
$ws.Range("K6").Value = 1.5
$ws.Range("L6").Formula = "=+D6*J6+1"
$ws.Range("M6").Formula = "=+ROUND(L6,0)"
$ws.Range("G6").Formula = "=+M6"

# --- Row 7 helper values + formulas ---
$ws.Range("B7").Value = 2
$ws.Range("J7").Value = 0.12
$ws.Range("K7").Value = 2.5
$ws.Range("L7").Formula = "=+D7*J7+1"
$ws.Range("M7").Formula = "=+ROUND(L7,0)"
$ws.Range("G7").Formula = "=+M7"

# --- Totals block ---
$ws.Range("G8").Formula = "=+SUM(G6:G7)"
$ws.Range("G9").Value = 0.07
$ws.Range("G10").Formula = "=+ROUND(G8*G9,0)"
$ws.Range("G11").Formula = "=+G8-G10"

# --- Column widths for the helper columns (best-fit, mirroring AutoFit) ---
$ws.Range("K1").ColumnWidth = 9.666666666666668
$ws.Range("M1").ColumnWidth = 11.5

# --- Row 5 height tweak ---
$ws.Range("A5").RowHeight = 19.5

# --- Data validations (drop-downs) ---
$dv1 = $ws.Range("E6:E7").Validation
$dv1.Delete()
$dv1.Add(3, 1, 1, '"Carta,Medio oficio,Oficio"')

$dv2 = $ws.Range("F6:F7").Validation
$dv2.Delete()
$dv2.Add(3, 1, 1, '"Blanco y negro,Colores"')

# --- View tweaks: zoom + selection ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("E6").Select()
